# Add 9 new exposure site rows at the top of the data (after the header row),
# shifting all existing rows down by 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows starting at row 2 (just below the header row)
$ws.Range("A2:A10").EntireRow.Insert()

# The freshly inserted rows inherit the header's bold/centered formatting;
# clear that so they match the plain formatting of the rest of the data rows.
$ws.Range("A2:D10").ClearFormats()

$newRows = @(
    @("1 Saintly Dr, Truganina VIC 3029", -37.807645, 144.746483, "Melton (C)"),
    @("300 Grattan St, Parkville VIC 3050", -37.798908, 144.956176, "Melbourne (C)"),
    @("176 Furlong Rd, St Albans VIC 3021", -37.759223, 144.816754, "Brimbank (C)"),
    @("185 Cooper St, Epping VIC 3076", -37.653023, 145.014685, "Whittlesea (C)"),
    @("10 Fletcher St, Essendon VIC 3040", -37.755206, 144.925092, "Moonee Valley (C)"),
    @("107-111 Andersons Creek Rd, Doncaster East VIC 3109", -37.766695, 145.170655, "Manningham (C)"),
    @("250 Waterdale Rd, Ivanhoe VIC 3079", -37.758211, 145.044162, "Banyule (C)"),
    @("27 Smith St, Healesville VIC 3777", -37.661037, 145.508141, "Yarra Ranges (S)"),
    @("17-19 Lalors Rd, Healesville VIC 3777", -37.66148, 145.516046, "Yarra Ranges (S)")
)

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
